# Adds changelog entries [1.13] and [1.14] to the "Versions" table (Sheet2 /
# sheet1.xml), grows the table / autofilter range to A1:C14, and turns on
# portrait page orientation - mirroring the authored OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft = -4131
$xlTop = -4160

# ---- Row 13 : [1.13] -------------------------------------------------
$a13 = $ws.Cells.Item(13, 1)
$a13.Value = "[1.13]"
$a13.HorizontalAlignment = $xlLeft
$a13.VerticalAlignment = $xlTop

$b13 = $ws.Cells.Item(13, 2)
$b13.Value = "[Printing]`n- implement postal cards - fix PC1 having extra page`n- use background worker for printing`n- fix a bug from last release that happens while generating envelops, they were not filled`n- implement using the same size for both envelops and postcards taken from UI"
$b13.HorizontalAlignment = $xlLeft
$b13.VerticalAlignment = $xlTop
$b13.WrapText = $true

$c13 = $ws.Cells.Item(13, 3)
$c13.Value = 43283
$c13.NumberFormat = "d-mmm-yy"
$c13.HorizontalAlignment = $xlLeft
$c13.VerticalAlignment = $xlTop

$ws.Rows.Item(13).RowHeight = 75

# ---- Row 14 : [1.14] -------------------------------------------------
$a14 = $ws.Cells.Item(14, 1)
$a14.Value = "[1.14]"
$a14.HorizontalAlignment = $xlLeft
$a14.VerticalAlignment = $xlTop

$b14 = $ws.Cells.Item(14, 2)
$b14.Value = "[Printing]`n- use the hotfix from Spire that fixes the footer issue`n- hide the other printing APIs options`n- adjust PC1 and PC2 to avoid overlapping Alberto's image when they get filled"
$b14.HorizontalAlignment = $xlLeft
$b14.VerticalAlignment = $xlTop
$b14.WrapText = $true

$c14 = $ws.Cells.Item(14, 3)
$c14.Value = 43304
$c14.NumberFormat = "d-mmm-yy"
$c14.HorizontalAlignment = $xlLeft
$c14.VerticalAlignment = $xlTop

$ws.Rows.Item(14).RowHeight = 60

# ---- grow the table / autofilter to include the two new rows ---------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C14"))

# ---- page setup: portrait orientation ---------------------------------
$ws.PageSetup.Orientation = 1

# ---- move selection to the newly added row, like the authored file ----
$ws.Range("B14").Select()

Write-Host "done"
